$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 114; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $v = $cell.Value()
    $cell.Value = $v.AddDays(1)
}
